# Updated cryptos list -- apply the new Price (D) / Volume(1h) (E) text values
# for each affected row, exactly as captured by the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.021.62'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '2.339.07'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''306.74'
$ws.Range("E5").Value = '  -1.47%  '
$ws.Range("D6").Value = '''101.02'
$ws.Range("E6").Value = '  -1.90%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -3.48%  '
$ws.Range("D10").Value = '''34.91'
$ws.Range("E10").Value = '  -2.26%  '
$ws.Range("D11").Value = '''52.58'
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("E12").Value = '  -2.01%  '
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("E14").Value = '  -2.41%  '
$ws.Range("D15").Value = '''15.80'
$ws.Range("E15").Value = '  +5.17%  '
$ws.Range("D16").Value = '2.321.60'
$ws.Range("E16").Value = '  -0.58%  '
$ws.Range("E17").Value = '  +2.42%  '
$ws.Range("D18").Value = '42.952.71'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("E19").Value = '  +0.91%  '
$ws.Range("D20").Value = '''11.73'
$ws.Range("E20").Value = '  -4.37%  '
$ws.Range("D21").Value = '0.0₃0910'
$ws.Range("E21").Value = '  -2.52%  '
$ws.Range("D23").Value = '''236.61'
$ws.Range("E23").Value = '  -2.04%  '
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("E25").Value = '  -1.83%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").Value = '''25.54'
$ws.Range("E27").Value = '  +2.66%  '
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("E29").Value = '  +5.23%  '
$ws.Range("D30").Value = '''35.72'
$ws.Range("E30").Value = '  -2.59%  '
$ws.Range("D31").Value = '''9.31'
$ws.Range("E31").Value = '  -3.61%  '
$ws.Range("D32").Value = '''163.98'
$ws.Range("E32").Value = '  -3.97%  '
$ws.Range("D33").Value = '''0.999'
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("E34").Value = '  -2.99%  '
$ws.Range("E35").Value = '  +7.69%  '
$ws.Range("D36").Value = '''17.55'
$ws.Range("E36").Value = '  -1.17%  '
$ws.Range("D37").Value = '''0.0727'
$ws.Range("E37").Value = '  -1.86%  '
$ws.Range("E38").Value = '  -4.16%  '
$ws.Range("E39").Value = '  -1.64%  '
$ws.Range("E40").Value = '  -4.83%  '
$ws.Range("E41").Value = '  -2.76%  '
$ws.Range("E42").Value = '  -2.60%  '
$ws.Range("D43").Value = '''2.59'
$ws.Range("E43").Value = '  +11.10%  '
$ws.Range("D44").Value = '2.026.44'
$ws.Range("E44").Value = '  +2.37%  '
$ws.Range("E45").Value = '  -2.36%  '
$ws.Range("D46").Value = '''18.87'
$ws.Range("E46").Value = '  -1.54%  '
$ws.Range("D47").Value = '''10.18'
$ws.Range("E47").Value = '  +1.61%  '
$ws.Range("E48").Value = '  -1.98%  '
$ws.Range("D49").Value = '''56.71'
$ws.Range("E49").Value = '  +2.30%  '
$ws.Range("D50").Value = '''2.90'
$ws.Range("E50").Value = '  -1.26%  '
$ws.Range("D51").Value = '2.564.72'
$ws.Range("E51").Value = '  +1.05%  '
